$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Seed the shared-string table in the exact order the strings first
# --- appear in the authored workbook (indices 13-20) ---
$ws.Range("A11").Value = "average "
$ws.Range("B11").Value = "short period"
$ws.Range("B12").Value = "spiral"
$ws.Range("A15").Value = "numerical"
$ws.Range("B16").Value = "phugoid"
$ws.Range("B13").Value = "dutch roll"
$ws.Range("B19").Value = "roll damping"
$ws.Range("A21").Value = "percentage difference"

# --- Column A width (closest achievable value to the authored 11.85546875) ---
$ws.Columns.Item(1).ColumnWidth = 11

# --- Row 11-13: averages of the measured (experimental) pairs ---
$ws.Range("C11").Formula = "=AVERAGE(C2:C3)"
$ws.Range("D11:F11").Formula = "=AVERAGE(D2:D3)"

$ws.Range("C12").Formula = "=AVERAGE(C5:C6)"
$ws.Range("D12:F12").Formula = "=AVERAGE(D5:D6)"

$ws.Range("C13").Formula = "=AVERAGE(C7:C8)"
$ws.Range("D13:F13").Formula = "=AVERAGE(D7:D8)"

# --- Row 15-19: numerical (reference) values ---
$ws.Range("B15").Value = "short period"
$ws.Range("C15").Value = 0.58
$ws.Range("D15").Value = -1.204
$ws.Range("E15").Value = 3.53
$ws.Range("F15").Value = 1.78

$ws.Range("C16").Value = 247.82
$ws.Range("D16").Value = -0.002797
$ws.Range("E16").Value = 46.47
$ws.Range("F16").Value = 0.1352

$ws.Range("B17").Value = "spiral"
$ws.Range("C17").Value = -57.71
$ws.Range("D17").Value = 0.01201

$ws.Range("B18").Value = "dutch roll"
$ws.Range("C18").Value = 2.56
$ws.Range("D18").Value = -0.271
$ws.Range("E18").Value = 3.3
$ws.Range("F18").Value = 1.905

$ws.Range("C19").Value = 0.18
$ws.Range("D19").Value = -3.865

# --- Row 21-25: percentage difference between averages/experiment and numerical values ---
$ws.Range("B21").Value = "short period"
$ws.Range("C21").Formula = "=(C11-C15)/C15*100"
$ws.Range("D21").Formula = "=(D11-D15)/D15*100"
$ws.Range("E21").Formula = "=(E11-E15)/E15*100"
$ws.Range("F21").Formula = "=(F11-F15)/F15*100"

$ws.Range("B22").Value = "phugoid"
$ws.Range("C22").Formula = "=(C4-C16)/C16*100"
$ws.Range("D22").Formula = "=(D4-D16)/D16*100"
$ws.Range("E22").Formula = "=(E4-E16)/E16*100"
$ws.Range("F22").Formula = "=(F4-F16)/F16*100"

$ws.Range("B23").Value = "spiral"
$ws.Range("C23").Formula = "=(C12-C17)/C17*100"
$ws.Range("D23").Formula = "=(D12-D17)/D17*100"
$ws.Range("E23").Formula = "=(E12-E17)/E17*100"
$ws.Range("F23").Formula = "=(F12-F17)/F17*100"

$ws.Range("B24").Value = "dutch roll"
$ws.Range("C24").Formula = "=(C13-C18)/C18*100"
$ws.Range("D24:F24").Formula = "=(D13-D18)/D18*100"

$ws.Range("B25").Value = "roll damping"
$ws.Range("C25").Formula = "=(C9-C19)/C19*100"
$ws.Range("D25").Formula = "=(D9-D19)/D19*100"
$ws.Range("E25:F25").Formula = "=(E9-E19)/E19*100"

# --- Window state / selection matching the committed view ---
$win = $wb.Windows.Item(1)
$win.WindowState = -4143
$ws.Range("D26").Select() | Out-Null
